{"js": "const replacements = [\n  [\"362\u00d75=1810\", \"230\u00d73=690\"],\n  [\"392\u00d74=1568\", \"641\u00d73=1923\"],\n  [\"792\u00d79=7128\", \"308\u00d79=2772\"],\n  [\"262\u00d76=1572\", \"291\u00d79=2619\"],\n  [\"409\u00d77=2863\", \"510\u00d72=1020\"],\n  [\"264\u00d73=792\", \"364\u00d76=2184\"],\n  [\"345\u00d79=3105\", \"131\u00d72=262\"],\n  [\"957\u00d78=7656\", \"231\u00d79=2079\"],\n  [\"184\u00d72=368\", \"542\u00d72=1084\"],\n  [\"698\u00d75=3490\", \"605\u00d79=5445\"],\n  [\"960\u00d78=7680\", \"127\u00d76=762\"],\n  [\"736\u00d72=1472\", \"472\u00d75=2360\"],\n  [\"992\u00d77=6944\", \"770\u00d79=6930\"],\n  [\"305\u00d75=1525\", \"413\u00d75=2065\"],\n  [\"302\u00d73=906\", \"693\u00d75=3465\"],\n  [\"827\u00d72=1654\", \"593\u00d78=4744\"],\n  [\"648\u00d76=3888\", \"495\u00d72=990\"],\n  [\"858\u00d72=1716\", \"909\u00d77=6363\"],\n  [\"917\u00d79=8253\", \"767\u00d75=3835\"],\n  [\"941\u00d78=7528\", \"221\u00d74=884\"],\n  [\"651\u00d77=4557\", \"942\u00d72=1884\"],\n  [\"280\u00d79=2520\", \"857\u00d73=2571\"],\n  [\"654\u00d79=5886\", \"461\u00d75=2305\"],\n  [\"524\u00d78=4192\", \"113\u00d78=904\"],\n  [\"272\u00d79=2448\", \"870\u00d79=7830\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"362\u00d75=1810\", \"230\u00d73=690\"),\n    @(\"392\u00d74=1568\", \"641\u00d73=1923\"),\n    @(\"792\u00d79=7128\", \"308\u00d79=2772\"),\n    @(\"262\u00d76=1572\", \"291\u00d79=2619\"),\n    @(\"409\u00d77=2863\", \"510\u00d72=1020\"),\n    @(\"264\u00d73=792\", \"364\u00d76=2184\"),\n    @(\"345\u00d79=3105\", \"131\u00d72=262\"),\n    @(\"957\u00d78=7656\", \"231\u00d79=2079\"),\n    @(\"184\u00d72=368\", \"542\u00d72=1084\"),\n    @(\"698\u00d75=3490\", \"605\u00d79=5445\"),\n    @(\"960\u00d78=7680\", \"127\u00d76=762\"),\n    @(\"736\u00d72=1472\", \"472\u00d75=2360\"),\n    @(\"992\u00d77=6944\", \"770\u00d79=6930\"),\n    @(\"305\u00d75=1525\", \"413\u00d75=2065\"),\n    @(\"302\u00d73=906\", \"693\u00d75=3465\"),\n    @(\"827\u00d72=1654\", \"593\u00d78=4744\"),\n    @(\"648\u00d76=3888\", \"495\u00d72=990\"),\n    @(\"858\u00d72=1716\", \"909\u00d77=6363\"),\n    @(\"917\u00d79=8253\", \"767\u00d75=3835\"),\n    @(\"941\u00d78=7528\", \"221\u00d74=884\"),\n    @(\"651\u00d77=4557\", \"942\u00d72=1884\"),\n    @(\"280\u00d79=2520\", \"857\u00d73=2571\"),\n    @(\"654\u00d79=5886\", \"461\u00d75=2305\"),\n    @(\"524\u00d78=4192\", \"113\u00d78=904\"),\n    @(\"272\u00d79=2448\", \"870\u00d79=7830\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n"}
